# Update run to cross multiple browsers: mark the chrome/edge/firefox rows
# on TESTDATA as executable, bump the VERSION values, and move the
# selection from G5 to F5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TESTDATA")
[void]$ws.Activate()

# Row 2 (chrome)
$ws.Range("C2").Value = 109
$ws.Range("D2").Value = "yes"

# Row 3 (edge)
$ws.Range("C3").Value = 110

# Row 4 (firefox)
$ws.Range("C4").Value = 111
$ws.Range("D4").Value = "yes"

[void]$ws.Range("F5").Select()
